$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header values B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: B2 and D2 get new values, C2 and E2 become empty
$ws.Range("B2").Value = 17.375431314037488
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 25.424661903179029
$ws.Range("E2").ClearContents()

# Row 3: update values B3:E3
$ws.Range("B3").Value = 14.471734435433772
$ws.Range("C3").Value = -10.616310651571711
$ws.Range("D3").Value = 18.216853546741003
$ws.Range("E3").Value = -15.401338659918059

# Update the selected range to reflect the new data extent (B1:E3)
[void]$ws.Range("B1:E3").Select()
